# Auto-generated Excel COM-interop script to apply scheduled price/profit refresh
# across the Lich_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets("ALC")
$ws.Range("H4").Value = 55.166668
$ws.Range("I4").Value = 59.4
$ws.Range("K4").Value = 59.4
$ws.Range("M4").Value = 54.6
$ws.Range("H5").Value = 126.22222
$ws.Range("I5").Value = 124.5
$ws.Range("K5").Value = 124.5
$ws.Range("M5").Value = -9.5
$ws.Range("H8").Value = 35.444443
$ws.Range("I8").Value = 28.714285
$ws.Range("K8").Value = 86.142855
$ws.Range("M8").Value = 52.857145
$ws.Range("H69").Value = 7354.304
$ws.Range("I69").Value = 7667.933
$ws.Range("J69").Value = 6766.25
$ws.Range("K69").Value = 23003.799
$ws.Range("L69").Value = 20298.75
$ws.Range("M69").Value = -22129.799
$ws.Range("N69").Value = -22046.75
$ws.Range("H72").Value = 7354.304
$ws.Range("I72").Value = 7667.933
$ws.Range("J72").Value = 6766.25
$ws.Range("K72").Value = 69011.397
$ws.Range("L72").Value = 60896.25
$ws.Range("M72").Value = -64643.397
$ws.Range("N72").Value = -69632.25
$ws.Range("H137").Value = 11382.294
$ws.Range("I137").Value = 14636.125
$ws.Range("J137").Value = 8490
$ws.Range("K137").Value = 43908.375
$ws.Range("L137").Value = 25470
$ws.Range("M137").Value = -41358.375
$ws.Range("N137").Value = -30570

# --- ARM ---
$ws = $wb.Worksheets("ARM")
$ws.Range("H3").Value = 10000
$ws.Range("J3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("N3").Value = -10230
$ws.Range("H61").Value = 5702.643
$ws.Range("I61").Value = 3354.111
$ws.Range("J61").Value = 9930
$ws.Range("K61").Value = 3354.111
$ws.Range("L61").Value = 9930
$ws.Range("M61").Value = -3142.111
$ws.Range("N61").Value = -10354
$ws.Range("H74").Value = 4181.906
$ws.Range("I74").Value = 4181.906
$ws.Range("K74").Value = 4181.906
$ws.Range("M74").Value = -3307.906
$ws.Range("H77").Value = 4181.906
$ws.Range("I77").Value = 4181.906
$ws.Range("K77").Value = 20909.53
$ws.Range("M77").Value = -16541.53
$ws.Range("H136").Value = 5702.643
$ws.Range("I136").Value = 3354.111
$ws.Range("J136").Value = 9930
$ws.Range("K136").Value = 10062.333
$ws.Range("L136").Value = 29790
$ws.Range("M136").Value = -7512.332999999999
$ws.Range("N136").Value = -34890

# --- BSM ---
$ws = $wb.Worksheets("BSM")
$ws.Range("H5").Value = 999.5
$ws.Range("I5").Value = 999.5
$ws.Range("K5").Value = 999.5
$ws.Range("M5").Value = -886.5
$ws.Range("H7").Value = 1001.5
$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = 110
$ws.Range("N7").Value = -2226
$ws.Range("H12").Value = 335.5
$ws.Range("I12").Value = 224.33333
$ws.Range("J12").Value = 446.66666
$ws.Range("K12").Value = 224.33333
$ws.Range("L12").Value = 446.66666
$ws.Range("M12").Value = -56.33332999999999
$ws.Range("N12").Value = -782.66666
$ws.Range("H86").Value = 1773.4
$ws.Range("I86").Value = 1604.6666
$ws.Range("K86").Value = 1604.6666
$ws.Range("M86").Value = -481.6666
$ws.Range("H89").Value = 1773.4
$ws.Range("I89").Value = 1604.6666
$ws.Range("K89").Value = 8023.333000000001
$ws.Range("M89").Value = -2407.333000000001

# --- CRP ---
$ws = $wb.Worksheets("CRP")
$ws.Range("H7").Value = 57
$ws.Range("I7").Value = 41.3
$ws.Range("K7").Value = 41.3
$ws.Range("M7").Value = 71.7
$ws.Range("H16").Value = 1663.5
$ws.Range("I16").Value = 1323.3334
$ws.Range("K16").Value = 1323.3334
$ws.Range("M16").Value = -1036.3334
$ws.Range("H31").Value = 45341.367
$ws.Range("I31").Value = 4158.2856
$ws.Range("J31").Value = 69364.836
$ws.Range("K31").Value = 4158.2856
$ws.Range("L31").Value = 69364.836
$ws.Range("M31").Value = -3863.2856
$ws.Range("N31").Value = -69954.836
$ws.Range("H34").Value = 45341.367
$ws.Range("I34").Value = 4158.2856
$ws.Range("J34").Value = 69364.836
$ws.Range("K34").Value = 4158.2856
$ws.Range("L34").Value = 69364.836
$ws.Range("M34").Value = -3956.2856
$ws.Range("N34").Value = -69768.836
$ws.Range("H99").Value = 2615.3333
$ws.Range("I99").Value = 2453.125
$ws.Range("J99").Value = 2745.1
$ws.Range("K99").Value = 2453.125
$ws.Range("L99").Value = 2745.1
$ws.Range("M99").Value = -955.125
$ws.Range("N99").Value = -5741.1
$ws.Range("H113").Value = 1663.5
$ws.Range("I113").Value = 1323.3334
$ws.Range("K113").Value = 1323.3334
$ws.Range("M113").Value = 846.6666
$ws.Range("H126").Value = 2615.3333
$ws.Range("I126").Value = 2453.125
$ws.Range("J126").Value = 2745.1
$ws.Range("K126").Value = 7359.375
$ws.Range("L126").Value = 8235.299999999999
$ws.Range("M126").Value = -4889.375
$ws.Range("N126").Value = -13175.3
$ws.Range("H134").Value = 5554.5
$ws.Range("I134").Value = 6236.5454
$ws.Range("J134").Value = 1803.25
$ws.Range("K134").Value = 18709.6362
$ws.Range("L134").Value = 5409.75
$ws.Range("M134").Value = -16174.6362
$ws.Range("N134").Value = -10479.75

# --- CUL ---
$ws = $wb.Worksheets("CUL")
$ws.Range("H4").Value = 10309033
$ws.Range("I4").Value = 10309033
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 30927099
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -30926987

# --- GSM ---
$ws = $wb.Worksheets("GSM")
$ws.Range("H2").Value = 82.90476
$ws.Range("I2").Value = 99.72727
$ws.Range("K2").Value = 99.72727
$ws.Range("M2").Value = 13.27273
$ws.Range("H3").Value = 1998
$ws.Range("I3").Value = 3250
$ws.Range("J3").Value = 1163.3334
$ws.Range("K3").Value = 3250
$ws.Range("L3").Value = 1163.3334
$ws.Range("M3").Value = -3134
$ws.Range("N3").Value = -1395.3334
$ws.Range("H4").Value = 3999
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 1799.5
$ws.Range("I5").Value = 1799.5
$ws.Range("K5").Value = 1799.5
$ws.Range("M5").Value = -1687.5
$ws.Range("H11").Value = 277084.5
$ws.Range("J11").Value = 302779.34
$ws.Range("L11").Value = 302779.34
$ws.Range("N11").Value = -303057.34
$ws.Range("H33").Value = 18900
$ws.Range("J33").Value = 18900
$ws.Range("L33").Value = 18900
$ws.Range("N33").Value = -19404

# --- LTW ---
$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 6612.778
$ws.Range("I7").Value = 6622.5674
$ws.Range("J7").Value = 6567.5
$ws.Range("K7").Value = 6622.5674
$ws.Range("L7").Value = 6567.5
$ws.Range("M7").Value = -6510.5674
$ws.Range("N7").Value = -6791.5
$ws.Range("H41").Value = 40007.75
$ws.Range("I41").Value = 40007.75
$ws.Range("K41").Value = 40007.75
$ws.Range("M41").Value = -39569.75
$ws.Range("H68").Value = 2661.0881
$ws.Range("J68").Value = 2498.5
$ws.Range("L68").Value = 2498.5
$ws.Range("N68").Value = -3996.5
$ws.Range("H71").Value = 2661.0881
$ws.Range("J71").Value = 2498.5
$ws.Range("L71").Value = 12492.5
$ws.Range("N71").Value = -19980.5
$ws.Range("H82").Value = 2084.5715
$ws.Range("I82").Value = 2326.2856
$ws.Range("J82").Value = 1842.8572
$ws.Range("K82").Value = 2326.2856
$ws.Range("L82").Value = 1842.8572
$ws.Range("M82").Value = -1965.2856
$ws.Range("N82").Value = -2564.8572
$ws.Range("H85").Value = 2084.5715
$ws.Range("I85").Value = 2326.2856
$ws.Range("J85").Value = 1842.8572
$ws.Range("K85").Value = 2326.2856
$ws.Range("L85").Value = 1842.8572
$ws.Range("M85").Value = -1078.2856
$ws.Range("N85").Value = -4338.8572
$ws.Range("H126").Value = 6612.778
$ws.Range("I126").Value = 6622.5674
$ws.Range("J126").Value = 6567.5
$ws.Range("K126").Value = 19867.7022
$ws.Range("L126").Value = 19702.5
$ws.Range("M126").Value = -17397.7022
$ws.Range("N126").Value = -24642.5

# --- WVR ---
$ws = $wb.Worksheets("WVR")
$ws.Range("H2").Value = 178946.33
$ws.Range("I2").Value = 265919.5
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 265919.5
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -265807.5
$ws.Range("N2").Value = -5224
$ws.Range("H126").Value = 5871.5
$ws.Range("I126").Value = 2498.6667
$ws.Range("J126").Value = 9244.333000000001
$ws.Range("K126").Value = 7496.000100000001
$ws.Range("L126").Value = 27732.999
$ws.Range("M126").Value = -5026.000100000001
$ws.Range("N126").Value = -32672.999
